$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("F2").Value = 12032
$ws.Range("F9").Value = 12160

$ws.Range("E10").Select()
